$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "September 19, 2025") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# ------------------------------------------------------------------
# 2. Split the mailing-address paragraph "3032 Lamory Pl, Santa Clara CA 95051"
#    (the one directly under the addressee name, NOT the one inside the
#    table further down) into two paragraphs:
#       "3032 Lamory Pl"
#       "Santa Clara, CA 95051"
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "3032 Lamory Pl, Santa Clara CA 95051") {
        $fullText = $p.Range.Text.TrimEnd()
        $commaIdx = $fullText.IndexOf(",")
        $splitPos = $p.Range.Start + $commaIdx

        # remove ", Santa Clara CA 95051" (keep the paragraph mark)
        $tailRange = $d.Range($splitPos, $p.Range.End - 1)
        $tailRange.Text = ""

        # insert a new paragraph break right after "3032 Lamory Pl"
        $breakRange = $d.Range($p.Range.End - 1, $p.Range.End - 1)
        $breakRange.InsertParagraphAfter()

        # fill in the new second line
        $newPara = $p.Next()
        $newPara.Range.Text = "Santa Clara, CA 95051"
        break
    }
}

# ------------------------------------------------------------------
# 3. Remove the now-superfluous blank "No Spacing" paragraph that used to
#    sit right after "...Board of Directors"
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Board of Directors*") {
        $next = $p.Next()
        if ($next.Range.Text.Trim() -eq "" -and $next.Style.NameLocal -eq "No Spacing") {
            $next.Range.Delete()
        }
        break
    }
}
